{"js": "// Map of old text -> new text for this edit (date header + 25 multiplication cells)\nconst replacements = [\n  [\"2025-09-15 Monday\", \"2025-09-16 Tuesday\"],\n  [\"70\u00d731=2170\", \"40\u00d721=840\"],\n  [\"48\u00d749=2352\", \"12\u00d785=1020\"],\n  [\"31\u00d752=1612\", \"11\u00d753=583\"],\n  [\"91\u00d783=7553\", \"44\u00d788=3872\"],\n  [\"32\u00d740=1280\", \"74\u00d711=814\"],\n  [\"50\u00d789=4450\", \"50\u00d745=2250\"],\n  [\"62\u00d757=3534\", \"71\u00d757=4047\"],\n  [\"25\u00d735=875\", \"34\u00d778=2652\"],\n  [\"53\u00d718=954\", \"33\u00d715=495\"],\n  [\"21\u00d731=651\", \"87\u00d776=6612\"],\n  [\"49\u00d723=1127\", \"65\u00d737=2405\"],\n  [\"89\u00d785=7565\", \"32\u00d792=2944\"],\n  [\"61\u00d780=4880\", \"55\u00d797=5335\"],\n  [\"13\u00d736=468\", \"25\u00d731=775\"],\n  [\"60\u00d711=660\", \"15\u00d759=885\"],\n  [\"95\u00d730=2850\", \"14\u00d796=1344\"],\n  [\"11\u00d789=979\", \"64\u00d741=2624\"],\n  [\"30\u00d734=1020\", \"60\u00d737=2220\"],\n  [\"25\u00d783=2075\", \"28\u00d786=2408\"],\n  [\"22\u00d782=1804\", \"96\u00d722=2112\"],\n  [\"90\u00d757=5130\", \"95\u00d732=3040\"],\n  [\"19\u00d737=703\", \"68\u00d797=6596\"],\n  [\"38\u00d723=874\", \"13\u00d737=481\"],\n  [\"12\u00d734=408\", \"57\u00d754=3078\"],\n  [\"57\u00d775=4275\", \"86\u00d793=7998\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date header and the 25 multiplication-table results.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-09-15 Monday\"; New = \"2025-09-16 Tuesday\" },\n    @{ Old = \"70\u00d731=2170\"; New = \"40\u00d721=840\" },\n    @{ Old = \"48\u00d749=2352\"; New = \"12\u00d785=1020\" },\n    @{ Old = \"31\u00d752=1612\"; New = \"11\u00d753=583\" },\n    @{ Old = \"91\u00d783=7553\"; New = \"44\u00d788=3872\" },\n    @{ Old = \"32\u00d740=1280\"; New = \"74\u00d711=814\" },\n    @{ Old = \"50\u00d789=4450\"; New = \"50\u00d745=2250\" },\n    @{ Old = \"62\u00d757=3534\"; New = \"71\u00d757=4047\" },\n    @{ Old = \"25\u00d735=875\"; New = \"34\u00d778=2652\" },\n    @{ Old = \"53\u00d718=954\"; New = \"33\u00d715=495\" },\n    @{ Old = \"21\u00d731=651\"; New = \"87\u00d776=6612\" },\n    @{ Old = \"49\u00d723=1127\"; New = \"65\u00d737=2405\" },\n    @{ Old = \"89\u00d785=7565\"; New = \"32\u00d792=2944\" },\n    @{ Old = \"61\u00d780=4880\"; New = \"55\u00d797=5335\" },\n    @{ Old = \"13\u00d736=468\"; New = \"25\u00d731=775\" },\n    @{ Old = \"60\u00d711=660\"; New = \"15\u00d759=885\" },\n    @{ Old = \"95\u00d730=2850\"; New = \"14\u00d796=1344\" },\n    @{ Old = \"11\u00d789=979\"; New = \"64\u00d741=2624\" },\n    @{ Old = \"30\u00d734=1020\"; New = \"60\u00d737=2220\" },\n    @{ Old = \"25\u00d783=2075\"; New = \"28\u00d786=2408\" },\n    @{ Old = \"22\u00d782=1804\"; New = \"96\u00d722=2112\" },\n    @{ Old = \"90\u00d757=5130\"; New = \"95\u00d732=3040\" },\n    @{ Old = \"19\u00d737=703\"; New = \"68\u00d797=6596\" },\n    @{ Old = \"38\u00d723=874\"; New = \"13\u00d737=481\" },\n    @{ Old = \"12\u00d734=408\"; New = \"57\u00d754=3078\" },\n    @{ Old = \"57\u00d775=4275\"; New = \"86\u00d793=7998\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $r.Old,        # FindText\n        $true,         # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap = wdFindContinue\n        $false,        # Format\n        $r.New,        # ReplaceWith\n        2              # Replace = wdReplaceAll\n    ) | Out-Null\n}\n"}
